$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 2 ("Jamaican Jerk Chicken" / "Chicken" / test values in G2:H2),
# shifting all subsequent rows up by one.
$ws.Rows.Item(2).Delete()

# Update the active selection to match the saved view state.
$ws.Range("A4").Select() | Out-Null
